# Apply the "new wfps regridding" data fix:
# - Rename "Flux, 1850" label rows to "Flux, 1860"
# - Update the "Global run" block values (Mean EF row, Effective EF row)
# - Fill in previously-empty Flux values for the "Global run" block
# - Fill in previously-empty Mean EF values for the "Global run, Arctic mask" block
# - Move the active selection to D13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel "Flux, 1850" -> "Flux, 1860" everywhere it appears ---
$ws.Range("B4").Value = "Flux, 1860"
$ws.Range("B9").Value = "Flux, 1860"
$ws.Range("B14").Value = "Flux, 1860"

# --- "Global run" block (rows 7-11) ---
# Mean EF row
$ws.Range("C7").Value = 1.1000000000000001
$ws.Range("D7").Value = 1.9

# Effective EF; 2022 row
$ws.Range("C8").Value = 4.3
$ws.Range("D8").Value = 7.2

# Flux, 1860 row
$ws.Range("C9").Value = 5.3

# Flux, 2020 row
$ws.Range("C10").Value = 14.1

# Flux, 2022 row
$ws.Range("C11").Value = 14.2

# --- "Global run, Arctic mask" block (rows 12-16) ---
# Mean EF row
$ws.Range("C12").Value = 0.14000000000000001
$ws.Range("D12").Value = 0.28999999999999998

# --- Update active selection ---
$ws.Range("D13").Select()
